$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Violet Hudson, poor"
$ws.Range("B2").Value = "Brooke Layton, good"
$ws.Range("C2").Value = "Niko Morris, poor"
$ws.Range("D2").Value = "Stanley Hirst, poor"
$ws.Range("E2").Value = "Caitlin Boyd, poor"
$ws.Range("F2").Value = "James Shilton, poor"
$ws.Range("G2").Value = "Lexi Green, poor"
$ws.Range("H2").Value = "James Calderon, poor"
$ws.Range("I2").Value = "Ava Lee, poor"
$ws.Range("J2").Value = "Esther Sido, excellent"

$ws.Range("B3").Value = "Katrina Petersone, good"
$ws.Range("D3").Value = "William Hunt, good"
$ws.Range("E3").Value = "Sophie Rayner, excellent"
$ws.Range("G3").Value = "Aarron Kelly, good"
$ws.Range("H3").Value = "Benjamin Finn, good"
$ws.Range("I3").Value = "Madison Taylor, good"
$ws.Range("J3").Value = "Benedict Hobday, good"

$ws.Range("B4").Value = "Alex Sentance, excellent"
$ws.Range("C4").Value = "Thomas Barrett, excellent"
$ws.Range("D4").Value = "Spencer Rowe, excellent"
$ws.Range("E4").Value = "Matthew Homan, excellent"
$ws.Range("F4").Value = "Benjamin Hillary, excellent"
$ws.Range("H4").Value = "Samuel Dixon, excellent"
